$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.306.92'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.681.02'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.03'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5533'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +8.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2701'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06479'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07554'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.539'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.672.90'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5803'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008431'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.99'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.340.59'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.922'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.91'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.34'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.221'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '146.75'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1315'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +9.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.900'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.79'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06329'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.390'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.81%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.590'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.574'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.665'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.036'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6181'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.400'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.714'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.235'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.112.92'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.90%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8709'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.68'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.831.13'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.34'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.173'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05274'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4291'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.062'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.27%  '
